$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A31").Value = "TieRod_tubeOffset (mm)"
$ws.Range("B31").Value = 20

$ws.Range("A31").NumberFormat = $ws.Range("A30").NumberFormat
$ws.Range("B31").NumberFormat = $ws.Range("B30").NumberFormat

# Reflect the updated selection/active cell recorded in the saved view state
$ws.Range("N8").Select() | Out-Null
